$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4701
$ws1.Range("F3").Value = 2553
$ws1.Range("F8").Value = 219
$ws1.Range("F10").Value = 1802
$ws1.Range("F12").Value = 4147
$ws1.Range("F14").Value = 286

# Sheet "全部类型" (all types) - update matching "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4701
$ws4.Range("F3").Value = 2553
$ws4.Range("F10").Value = 219
$ws4.Range("F14").Value = 1802
$ws4.Range("F16").Value = 4147
$ws4.Range("F18").Value = 286

$wb.Save()
